# Appends the newest Adafruit IO reading as row 77 (A1:F76 -> A1:F77).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
if ($row -lt 2) { $row = 2 }

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# Column C holds numeric-looking readings that must stay text (matches the
# existing inline-string "25" values used throughout the sheet), so force
# text formatting before assigning the value to avoid Excel auto-converting
# it to a number.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "25"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
